# Refresh the cryptos list snapshot (price / 1h-volume-change columns,
# plus the Stellar/Cosmos rows trading rank places) to match the latest
# coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.312.40'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.559.83'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '''0.490'
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''44.31'
$ws.Range('E8').Value = '  -4.55%  '
$ws.Range('D9').Value = '''23.57'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').Value = '''0.0894'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = '1.782.00'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '1.561.22'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '28.302.27'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '''0.511'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '''60.92'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').Value = '''227.40'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('D32').Value = '''1.05'
$ws.Range('E32').Value = '  -4.79%  '
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('D35').Value = '1.375.79'
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('D42').Value = '''1.93'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '''0.0471'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '''5.32'
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('D47').Value = '''62.04'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D49').Value = '1.695.12'
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').Value = '''85.23'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('E51').Value = '  -2.12%  '

# Stellar (row 28) and Cosmos (row 29) swapped ranking positions.
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''6.32'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '''0.103'
$ws.Range('E29').Value = '  -0.55%  '
